# Absenzenlisten-Templates 2016/2017 minimal ueberarbeitet
# (zu grosse Schrift in einzelnen Zellen korrigiert)
#
# Several single-letter "F" cells in the attendance table were still using
# the document's default run size (11pt / w:sz 22) instead of the 10pt
# (w:sz 20 / w:szCs 20) used everywhere else in the table. This script
# walks the table and, for every "F" cell that is still at the default
# size, sets it (and the immediately following cell, whose paragraph mark
# was missing the matching size) to 10pt.

$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    foreach ($row in $t.Rows) {
        $cells = $row.Cells
        $nCells = $cells.Count
        for ($ci = 1; $ci -le $nCells; $ci++) {
            $cell = $cells.Item($ci)
            $txt = $cell.Range.Text

            # Cell.Range.Text includes the trailing cell-mark (and possibly
            # a paragraph mark), so a lone "F" cell has length <= ~4 and
            # starts with "F".
            if ($txt.Length -ge 1 -and $txt.Length -le 4 -and $txt.Substring(0, 1) -eq "F") {

                if ($cell.Range.Font.Size -eq 11) {
                    # Fix the oversized "F" cell itself (run + paragraph mark).
                    $cell.Range.Font.Size = 10
                    $cell.Range.Font.SizeBi = 10

                    # The next cell in the row has the correct run size
                    # already, but is missing the matching paragraph-mark
                    # (w:pPr/w:rPr) size - apply the same, idempotent for
                    # the run, additive for the mark.
                    if ($ci -lt $nCells) {
                        $nextCell = $cells.Item($ci + 1)
                        $nextCell.Range.Font.Size = 10
                        $nextCell.Range.Font.SizeBi = 10
                    }
                }
            }
        }
    }
}
